$d = $word.ActiveDocument

$d.Content.Find.Execute("620÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "821÷9=", 2) | Out-Null
$d.Content.Find.Execute("765÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "863÷4=", 2) | Out-Null
$d.Content.Find.Execute("755÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "828÷9=", 2) | Out-Null
$d.Content.Find.Execute("478÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "424÷7=", 2) | Out-Null
$d.Content.Find.Execute("872÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "240÷2=", 2) | Out-Null
$d.Content.Find.Execute("357÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "925÷2=", 2) | Out-Null
$d.Content.Find.Execute("612÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "241÷3=", 2) | Out-Null
$d.Content.Find.Execute("989÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "613÷2=", 2) | Out-Null
$d.Content.Find.Execute("860÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "272÷2=", 2) | Out-Null
$d.Content.Find.Execute("549÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "924÷6=", 2) | Out-Null
$d.Content.Find.Execute("401÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "567÷5=", 2) | Out-Null
$d.Content.Find.Execute("250÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "833÷7=", 2) | Out-Null
$d.Content.Find.Execute("715÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "105÷2=", 2) | Out-Null
$d.Content.Find.Execute("346÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "704÷9=", 2) | Out-Null
$d.Content.Find.Execute("906÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "978÷3=", 2) | Out-Null
$d.Content.Find.Execute("313÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "688÷8=", 2) | Out-Null
$d.Content.Find.Execute("291÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "287÷6=", 2) | Out-Null
$d.Content.Find.Execute("566÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "259÷2=", 2) | Out-Null
$d.Content.Find.Execute("408÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "949÷9=", 2) | Out-Null
$d.Content.Find.Execute("356÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "570÷4=", 2) | Out-Null
$d.Content.Find.Execute("530÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "990÷8=", 2) | Out-Null
$d.Content.Find.Execute("648÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "441÷6=", 2) | Out-Null
$d.Content.Find.Execute("139÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "461÷8=", 2) | Out-Null
$d.Content.Find.Execute("670÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "587÷9=", 2) | Out-Null
$d.Content.Find.Execute("416÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "446÷7=", 2) | Out-Null
